# "Ue + moyenne par competence"
# Add the UE ("1/1") label and the per-competence averages for student 1
# (row 2) in columns M:P.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "1/1"
$ws.Range("N2").Value = 17
$ws.Range("O2").Value = 17.5
$ws.Range("P2").Value = 15.5
